$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last updated" timestamp banner in A1
$ws.Range("A1").Value = "Last updated: 2025-07-15 14:22:10"

# Swap the CommittedNotShip (C21) and UncommittedOrders (D21) values for
# PO 4516351202_TIPI
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 5
